{"js": "// Replace the date line and all of the multiplication problems in the\n// table with the new values from the target revision. Every old value is\n// unique in the document, so a simple exact-text search/replace per pair\n// is safe and unambiguous.\nconst replacements = [\n  [\"2026-02-28 Saturday\", \"2026-03-01 Sunday\"],\n  [\"913\u00d72=\", \"441\u00d74=\"],\n  [\"661\u00d73=\", \"380\u00d72=\"],\n  [\"703\u00d73=\", \"818\u00d75=\"],\n  [\"358\u00d75=\", \"300\u00d76=\"],\n  [\"726\u00d79=\", \"984\u00d79=\"],\n  [\"124\u00d76=\", \"442\u00d77=\"],\n  [\"554\u00d78=\", \"472\u00d75=\"],\n  [\"683\u00d79=\", \"217\u00d75=\"],\n  [\"506\u00d79=\", \"879\u00d75=\"],\n  [\"785\u00d79=\", \"969\u00d75=\"],\n  [\"803\u00d78=\", \"349\u00d77=\"],\n  [\"180\u00d77=\", \"842\u00d75=\"],\n  [\"231\u00d73=\", \"176\u00d77=\"],\n  [\"152\u00d72=\", \"846\u00d72=\"],\n  [\"946\u00d76=\", \"823\u00d77=\"],\n  [\"392\u00d76=\", \"509\u00d76=\"],\n  [\"636\u00d79=\", \"546\u00d76=\"],\n  [\"611\u00d77=\", \"179\u00d79=\"],\n  [\"659\u00d74=\", \"275\u00d75=\"],\n  [\"667\u00d74=\", \"780\u00d78=\"],\n  [\"181\u00d77=\", \"470\u00d78=\"],\n  [\"961\u00d74=\", \"811\u00d77=\"],\n  [\"650\u00d74=\", \"240\u00d78=\"],\n  [\"328\u00d79=\", \"383\u00d75=\"],\n  [\"592\u00d72=\", \"947\u00d73=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and all of the multiplication problems in the\n# table with the new values from the target revision. Every old value is\n# unique in the document, so a simple Find/Replace per pair is safe and\n# unambiguous.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"2026-02-28 Saturday\", \"2026-03-01 Sunday\"),\n  @(\"913\u00d72=\", \"441\u00d74=\"),\n  @(\"661\u00d73=\", \"380\u00d72=\"),\n  @(\"703\u00d73=\", \"818\u00d75=\"),\n  @(\"358\u00d75=\", \"300\u00d76=\"),\n  @(\"726\u00d79=\", \"984\u00d79=\"),\n  @(\"124\u00d76=\", \"442\u00d77=\"),\n  @(\"554\u00d78=\", \"472\u00d75=\"),\n  @(\"683\u00d79=\", \"217\u00d75=\"),\n  @(\"506\u00d79=\", \"879\u00d75=\"),\n  @(\"785\u00d79=\", \"969\u00d75=\"),\n  @(\"803\u00d78=\", \"349\u00d77=\"),\n  @(\"180\u00d77=\", \"842\u00d75=\"),\n  @(\"231\u00d73=\", \"176\u00d77=\"),\n  @(\"152\u00d72=\", \"846\u00d72=\"),\n  @(\"946\u00d76=\", \"823\u00d77=\"),\n  @(\"392\u00d76=\", \"509\u00d76=\"),\n  @(\"636\u00d79=\", \"546\u00d76=\"),\n  @(\"611\u00d77=\", \"179\u00d79=\"),\n  @(\"659\u00d74=\", \"275\u00d75=\"),\n  @(\"667\u00d74=\", \"780\u00d78=\"),\n  @(\"181\u00d77=\", \"470\u00d78=\"),\n  @(\"961\u00d74=\", \"811\u00d77=\"),\n  @(\"650\u00d74=\", \"240\u00d78=\"),\n  @(\"328\u00d79=\", \"383\u00d75=\"),\n  @(\"592\u00d72=\", \"947\u00d73=\")\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $find = $d.Content.Find\n  $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
